$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 2.2

# Row 4 updates
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
